# "Copy of Form nghi phep.xlsx" - add an "approval date" column.
#
# The sheet was a 15-row blank template; it gets trimmed down to just the
# header row + a single blank data row, and a new "Ngay duyet don" (Approval
# date) column is inserted right after "Ngay tao don", pushing "Ly do" and
# "Trang thai" one column over (F/G/H).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop the extra blank template rows; keep header (row 1) + one data row (row 2) ---
$ws.Rows("3:15").Delete()

# --- Shift headers right: H <- G <- F, then put the new header into F ---
$ws.Range("H1").Value = $ws.Range("G1").Value2
$ws.Range("G1").Value = $ws.Range("F1").Value2
$ws.Range("F1").Value = "Ngày duyệt đơn"

# --- New F1 header look: same yellow fill as the rest of the header row,
#     centered, but with no cell border ---
$ws.Range("F1").Interior.Color = 65535
$ws.Range("F1").Borders.LineStyle = 0
$ws.Range("F1").HorizontalAlignment = -4108
$ws.Range("F1").VerticalAlignment = -4108

# --- F2 is no longer the "creation date" column, so give it the plain bordered
#     text look (like A2/B2) instead of the date-formatted one ---
$ws.Range("A2").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Column F can be narrower now that it just holds a date ---
$ws.Columns("F").ColumnWidth = 21.5

# --- Reset the view: scroll back to A1 and move the active selection ---
$ws.Range("A1").Select()
$ws.Range("G4").Select()
